$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.85
$ws.Range("I2").Value = 5.1
$ws.Range("P2").Value = 1.89
$ws.Range("AA2").Value = 130
$ws.Range("AI2").Value = 70
$ws.Range("AO2").Value = 70
$ws.Range("G3").Value = 2.58
$ws.Range("I3").Value = 3.6
$ws.Range("P3").Value = 1.87
$ws.Range("Q3").Value = 1.92
$ws.Range("R3").Value = 1.33
$ws.Range("T3").Value = 1.74
$ws.Range("V3").Value = 1.38
$ws.Range("W3").Value = 1.63
$ws.Range("AD3").Value = 17
$ws.Range("F4").Value = 1.37
$ws.Range("G4").Value = 1.46
$ws.Range("K4").Value = 7.2
$ws.Range("L4").Value = 1.22
$ws.Range("O4").Value = 1.16
$ws.Range("P4").Value = 2.64
$ws.Range("Q4").Value = 1.5
$ws.Range("S4").Value = 2.12
$ws.Range("T4").Value = 1.64
$ws.Range("U4").Value = 1.92
$ws.Range("W4").Value = 3.15
$ws.Range("X4").Value = 32
$ws.Range("Y4").Value = 38
$ws.Range("Z4").Value = 100
$ws.Range("AD4").Value = 1000
$ws.Range("AG4").Value = 11.5
$ws.Range("K5").Value = 7.4
$ws.Range("N5").Value = 5.8
$ws.Range("Q5").Value = 1.52
$ws.Range("R5").Value = 1.65
$ws.Range("T5").Value = 2.04
$ws.Range("F6").Value = 2.46
$ws.Range("N6").Value = 3.15
$ws.Range("P6").Value = 1.74
$ws.Range("T6").Value = 1.89
$ws.Range("V6").Value = 1.4
$ws.Range("F7").Value = 1.46
$ws.Range("G7").Value = 1.6
$ws.Range("H7").Value = 5.5
$ws.Range("I7").Value = 8.800000000000001
$ws.Range("N7").Value = 4.9
$ws.Range("Q7").Value = 1.5
$ws.Range("R7").Value = 1.63
$ws.Range("U7").Value = 1.99
$ws.Range("V7").Value = 1.13
$ws.Range("W7").Value = 2.68
$ws.Range("Z7").Value = 75
$ws.Range("AE7").Value = 100
$ws.Range("AN7").Value = 6.6
$ws.Range("AO7").Value = 90
$ws.Range("L8").Value = 1.25
$ws.Range("I9").Value = 1.95
$ws.Range("J9").Value = 3.45
$ws.Range("K9").Value = 4
$ws.Range("N9").Value = 3.3
$ws.Range("V9").Value = 2.04
$ws.Range("F10").Value = 7.4
$ws.Range("J10").Value = 5.3
$ws.Range("F11").Value = 1.54
$ws.Range("G11").Value = 1.56
$ws.Range("I11").Value = 9.4
$ws.Range("J11").Value = 4.2
$ws.Range("K11").Value = 4.3
$ws.Range("P11").Value = 1.67
$ws.Range("Q11").Value = 2.2
$ws.Range("V11").Value = 1.12
$ws.Range("W11").Value = 2.78
$ws.Range("AH11").Value = 34
$ws.Range("AI11").Value = 200
$ws.Range("AM11").Value = 290
$ws.Range("F13").Value = 1.46
$ws.Range("L13").Value = 1.34
$ws.Range("O13").Value = 1.25
$ws.Range("Q13").Value = 1.73
$ws.Range("S13").Value = 2.84
$ws.Range("L15").Value = 1.45
$ws.Range("M15").Value = 1.1
$ws.Range("N15").Value = 2.82
$ws.Range("O15").Value = 1.44
$ws.Range("P15").Value = 1.64
$ws.Range("Q15").Value = 2.28
$ws.Range("R15").Value = 1.22
$ws.Range("S15").Value = 4.5
$ws.Range("T15").Value = 1.93
$ws.Range("G16").Value = 1.84
$ws.Range("I16").Value = 8.4
$ws.Range("J16").Value = 3.45
$ws.Range("O16").Value = 1.35
$ws.Range("F17").Value = 3.95
$ws.Range("G17").Value = 4.8
$ws.Range("H17").Value = 2.06
$ws.Range("I17").Value = 2.36
$ws.Range("K17").Value = 3.85
$ws.Range("L17").Value = 1.51
$ws.Range("N17").Value = 2.76
$ws.Range("S17").Value = 4.1
$ws.Range("T17").Value = 1.86
$ws.Range("U17").Value = 1.68
$ws.Range("V17").Value = 1.76
$ws.Range("W17").Value = 1.27
$ws.Range("I18").Value = 2.16
$ws.Range("J18").Value = 3.2
$ws.Range("Q18").Value = 2.06
$ws.Range("S18").Value = 4.2
$ws.Range("T18").Value = 1.8
$ws.Range("U18").Value = 1.73
$ws.Range("V18").Value = 1.86
$ws.Range("H19").Value = 4
$ws.Range("R19").Value = 1.26
$ws.Range("H20").Value = 4.1
$ws.Range("J20").Value = 3.75
$ws.Range("P20").Value = 2.04
$ws.Range("Q20").Value = 1.79
$ws.Range("F21").Value = 1.86
$ws.Range("G21").Value = 1.87
$ws.Range("J21").Value = 4.2
$ws.Range("K21").Value = 4.3
$ws.Range("W21").Value = 2.14
$ws.Range("AJ21").Value = 21
$ws.Range("AK21").Value = 16
$ws.Range("AM21").Value = 60
$ws.Range("H22").Value = 2.84
$ws.Range("H23").Value = 19
$ws.Range("I23").Value = 21
$ws.Range("J23").Value = 7.6
$ws.Range("K23").Value = 8.199999999999999
$ws.Range("N23").Value = 5
$ws.Range("Q23").Value = 1.54
$ws.Range("R23").Value = 1.55
$ws.Range("S23").Value = 2.52
$ws.Range("T23").Value = 2.52
$ws.Range("U23").Value = 1.54
$ws.Range("V23").Value = 1.05
$ws.Range("AE23").Value = 560
$ws.Range("AI23").Value = 380
$ws.Range("AL23").Value = 60
$ws.Range("F24").Value = 3.3
$ws.Range("G24").Value = 3.5
$ws.Range("I24").Value = 2.18
$ws.Range("J24").Value = 4.2
$ws.Range("V24").Value = 1.84
$ws.Range("W24").Value = 1.4
$ws.Range("AJ24").Value = 70
$ws.Range("AK24").Value = 38
